$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- New column J: hours_worked ---
$ws.Range("J1").Value = "hours_worked"
# Match the bold/bordered header formatting already used by A1:I1.
$ws.Range("I1").Copy()
$ws.Range("J1").PasteSpecial(-4122)
$ws.Range("J2").Value = "2:02:00"
$ws.Range("J3").Value = "2:53:24"
$ws.Range("J4").Value = "0:00:48"
$ws.Range("J5").Value = "1:03:22"
$ws.Range("J6").Value = "0:00:29"

# --- New rows 7-14 ---
# Columns A (idtimeentry) and H (employeeid) hold plain integer strings which
# Excel would otherwise auto-convert to numbers, so force them to text by
# temporarily applying a text number format to those columns before writing,
# then restoring the default "Normal" style afterwards (matching the source
# workbook, where every data cell -- other than the header row -- carries no
# explicit style).
$ws.Range("A7:A14").NumberFormat = "@"
$ws.Range("H7:H14").NumberFormat = "@"

$ws.Range("A7").Value = "26"
$ws.Range("B7").Value = "2024-05-31 00:00:00"
$ws.Range("C7").Value = "18:32:53"
$ws.Range("D7").Value = "0:00:00"
$ws.Range("E7").Value = "2024-05-31 18:33:02"
$ws.Range("F7").Value = "2024-05-31 18:33:02"
$ws.Range("G7").Value = "None"
$ws.Range("H7").Value = "3"
$ws.Range("I7").Value = "alexa rodrig"
$ws.Range("J7").Value = "-1 day, 5:27:07"

$ws.Range("A8").Value = "29"
$ws.Range("B8").Value = "2024-05-31 00:00:00"
$ws.Range("C8").Value = "18:48:00"
$ws.Range("D8").Value = "18:53:44"
$ws.Range("E8").Value = "2024-05-31 18:48:00"
$ws.Range("F8").Value = "2024-05-31 22:53:44"
$ws.Range("G8").Value = "None"
$ws.Range("H8").Value = "5"
$ws.Range("I8").Value = "erit gridnev"
$ws.Range("J8").Value = "0:05:44"

$ws.Range("A9").Value = "30"
$ws.Range("B9").Value = "2024-05-31 00:00:00"
$ws.Range("C9").Value = "22:23:26"
$ws.Range("D9").Value = "23:53:33"
$ws.Range("E9").Value = "2024-05-31 22:23:26"
$ws.Range("F9").Value = "2024-06-01 03:53:33"
$ws.Range("G9").Value = "None"
$ws.Range("H9").Value = "4"
$ws.Range("I9").Value = "Avi Semah"
$ws.Range("J9").Value = "1:30:07"

$ws.Range("A10").Value = "31"
$ws.Range("B10").Value = "2024-05-31 00:00:00"
$ws.Range("C10").Value = "None"
$ws.Range("D10").Value = "23:21:00"
$ws.Range("E10").Value = "2024-05-31 23:21:00"
$ws.Range("F10").Value = "2024-05-31 23:21:00"
$ws.Range("G10").Value = "None"
$ws.Range("H10").Value = "4"
$ws.Range("I10").Value = "Avi Semah"
$ws.Range("J10").Value = "None"

$ws.Range("A11").Value = "32"
$ws.Range("B11").Value = "2024-05-31 00:00:00"
$ws.Range("C11").Value = "23:22:11"
$ws.Range("D11").Value = "23:22:47"
$ws.Range("E11").Value = "2024-05-31 23:22:11"
$ws.Range("F11").Value = "2024-05-31 23:22:47"
$ws.Range("G11").Value = "None"
$ws.Range("H11").Value = "4"
$ws.Range("I11").Value = "Avi Semah"
$ws.Range("J11").Value = "0:00:36"

$ws.Range("A12").Value = "33"
$ws.Range("B12").Value = "2024-05-31 00:00:00"
$ws.Range("C12").Value = "23:46:30"
$ws.Range("D12").Value = "23:50:27"
$ws.Range("E12").Value = "2024-06-01 03:46:30"
$ws.Range("F12").Value = "2024-05-31 23:50:27"
$ws.Range("G12").Value = "None"
$ws.Range("H12").Value = "4"
$ws.Range("I12").Value = "Avi Semah"
$ws.Range("J12").Value = "0:03:57"

$ws.Range("A13").Value = "34"
$ws.Range("B13").Value = "2024-05-31 00:00:00"
$ws.Range("C13").Value = "23:51:23"
$ws.Range("D13").Value = "23:58:08"
$ws.Range("E13").Value = "2024-05-31 23:51:23"
$ws.Range("F13").Value = "2024-05-31 23:58:08"
$ws.Range("G13").Value = "None"
$ws.Range("H13").Value = "4"
$ws.Range("I13").Value = "Avi Semah"
$ws.Range("J13").Value = "0:06:45"

$ws.Range("A14").Value = "35"
$ws.Range("B14").Value = "2024-06-01 00:00:00"
$ws.Range("C14").Value = "0:08:21"
$ws.Range("D14").Value = "0:09:10"
$ws.Range("E14").Value = "2024-06-01 04:08:22"
$ws.Range("F14").Value = "2024-06-01 04:09:11"
$ws.Range("G14").Value = "None"
$ws.Range("H14").Value = "4"
$ws.Range("I14").Value = "Avi Semah"
$ws.Range("J14").Value = "0:00:49"

# Restore the default (unstyled) appearance for the columns we forced to text.
$ws.Range("A7:A14").Style = "Normal"
$ws.Range("H7:H14").Style = "Normal"
